# Apply updated error-model values on the "Comparacion" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comparacion")

$ws.Range("E3").Value  = 3.7322535
$ws.Range("H3").Value  = 3.35198412
$ws.Range("J3").Value  = 1.64153286
$ws.Range("M3").Value  = 0.07530045

$ws.Range("E4").Value  = 3.57582647
$ws.Range("H4").Value  = 3.59696234

$ws.Range("E5").Value  = 3.59021725
$ws.Range("H5").Value  = 3.83023835

$ws.Range("E6").Value  = 4.16163355
$ws.Range("H6").Value  = 4.53630881

$ws.Range("E7").Value  = 4.62342302
$ws.Range("H7").Value  = 4.91372624

$ws.Range("E8").Value  = 5.35126635
$ws.Range("H8").Value  = 5.45871984

$ws.Range("E9").Value  = 6.26859618
$ws.Range("H9").Value  = 6.13968788

$ws.Range("E10").Value = 6.79831853
$ws.Range("H10").Value = 6.54812268

$ws.Range("E11").Value = 7.99987812
$ws.Range("H11").Value = 7.55303439

$ws.Range("E12").Value = 9.89671435
$ws.Range("H12").Value = 9.49614301

$ws.Range("E13").Value = 10.9714568
$ws.Range("H13").Value = 10.89795819

$ws.Range("E14").Value = 12.13041548
$ws.Range("H14").Value = 12.78531975
